# "Varios cambios sin mucha importancia"
# Update the "Caja SAC" row (row 3) of Hoja1 with the current list of
# providers/payment collection points. Google / Amazon / Efectivo / Sumset
# are replaced by Servientrega / dimonex / banco caja social / epa, while
# bancolombia / davivienda / banco bbva are kept as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = "Servientrega"
$ws.Range("C3").Value = "bancolombia"
$ws.Range("D3").Value = "davivienda"
$ws.Range("E3").Value = "dimonex"
$ws.Range("F3").Value = "banco caja social"
$ws.Range("G3").Value = "epa"
$ws.Range("H3").Value = "banco bbva"

# The user's selection ended up on H3 before saving.
$ws.Range("H3").Select()
